$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.347.01"
$ws.Range("E2").Value = "  +1.24%  "

$ws.Range("D3").Value = "1.622.59"
$ws.Range("E3").Value = "  +1.62%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.70"
$ws.Range("E10").Value = "  +2.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  +0.64%  "

$ws.Range("D12").Value = "1.849.63"
$ws.Range("E12").Value = "  +1.71%  "

$ws.Range("D13").Value = "1.632.03"
$ws.Range("E13").Value = "  +2.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.03"
$ws.Range("E14").Value = "  +0.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.518"
$ws.Range("E15").Value = "  +0.24%  "

$ws.Range("D16").Value = "26.356.44"
$ws.Range("E16").Value = "  +1.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.44"
$ws.Range("E17").Value = "  +2.62%  "

$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.50"
$ws.Range("E20").Value = "  -0.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.26"
$ws.Range("E21").Value = "  +0.16%  "

$ws.Range("E22").Value = "  +0.55%  "

$ws.Range("E23").Value = "  +0.41%  "

$ws.Range("E24").Value = "  -2.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.68"
$ws.Range("E25").Value = "  +0.65%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.120"
$ws.Range("E27").Value = "  -1.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.20"
$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.57"
$ws.Range("E29").Value = "  +0.69%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0520"
$ws.Range("E30").Value = "  +9.16%  "

$ws.Range("E31").Value = "  +0.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.18"
$ws.Range("E32").Value = "  +1.45%  "

$ws.Range("E33").Value = "  +0.51%  "

$ws.Range("E34").Value = "  +0.42%  "

$ws.Range("E35").Value = "  +2.32%  "

$ws.Range("D36").Value = "1.157.10"
$ws.Range("E36").Value = "  +2.32%  "

$ws.Range("E37").Value = "  +0.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.805"
$ws.Range("E38").Value = "  +1.19%  "

$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("E41").Value = "  +0.97%  "

$ws.Range("E42").Value = "  +4.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.784"
$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("D44").Value = "1.760.54"
$ws.Range("E44").Value = "  +1.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.66"
$ws.Range("E45").Value = "  +0.66%  "

$ws.Range("E46").Value = "  +10.04%  "

$ws.Range("E47").Value = "  +0.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.83"
$ws.Range("E48").Value = "  -0.30%  "

$ws.Range("E49").Value = "  +0.79%  "

$ws.Range("E50").Value = "  +1.03%  "

$ws.Range("E51").Value = "  -0.43%  "
